{"js": "// Append 3 more paragraphs identical to the existing ones\n// (\"Documento de Evidencias - DemoBlaze\", font size 18pt / w:sz=36)\n// to the end of the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst text = \"Documento de Evidencias - DemoBlaze\";\nconst fontSizePt = 18; // w:sz is in half-points -> 36 == 18pt\n\nlet last = paragraphs.items[paragraphs.items.length - 1];\n\nfor (let i = 0; i < 3; i++) {\n  const newPara = last.insertParagraph(text, \"After\");\n  newPara.font.size = fontSizePt;\n  last = newPara;\n}\n\nawait context.sync();\n", "ps1": "# Append 3 more paragraphs identical to the existing ones\n# (\"Documento de Evidencias - DemoBlaze\", font size 18pt / w:sz=36)\n# to the end of the document body.\n\n$d = $word.ActiveDocument\n$text = \"Documento de Evidencias - DemoBlaze\"\n\nfor ($i = 0; $i -lt 3; $i++) {\n    $last = $d.Paragraphs.Last\n    $last.Range.InsertParagraphAfter()\n    $newLast = $d.Paragraphs.Last\n    $newLast.Range.Text = $text\n}\n"}
